$d = $word.ActiveDocument

# The "fig-map" figure (an inline picture pulled from a now-dead Notion
# download link) and its caption paragraph ("La Palma is one of the west
# most islands in the Volcanic Archipelago of the Canary Islands
# (?@fig-map).") were dropped from the source. Remove that whole
# paragraph, identified via the "fig-map" bookmark that marks the
# picture inside it.
$bm = $d.Bookmarks("fig-map")
$p = $bm.Range.Paragraphs.First
$p.Range.Delete()
